# Rota.xlsx - "Criação da página ramais.tml"
#
# The schedule entry in C2 ("CARLOS") is corrected to "vicente", and the
# unused, never-formatted trailing columns (G:I) and trailing blank rows
# (6:7) plus the leftover style on F8 are cleaned out, shrinking the
# sheet's used range down to A1:F8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the collaborator name on the Saturday morning shift.
$ws.Range("C2").Value = "vicente"

# Clear out the stray formatted-but-empty columns G:I (rows 1-8) ...
$ws.Range("G1:I8").Clear()

# ... and the stray formatted-but-empty cells F6:F8.
$ws.Range("F6:F8").Clear()

# Move the active selection to D9 (matches the saved view state).
$ws.Range("D9").Select() | Out-Null
